$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 9.796688323618699
$ws.Cells.Item(2, 3).Value = -8.419353553388079
$ws.Cells.Item(2, 4).Value = -0.6232160131165674
$ws.Cells.Item(2, 5).Value = 0.719781571001505
$ws.Cells.Item(2, 6).Value = -1.645709516164073
$ws.Cells.Item(2, 7).Value = -0.9785742016169946
$ws.Cells.Item(2, 8).Value = -0.20769435928658
$ws.Cells.Item(2, 9).Value = -0.2996160391110768
$ws.Cells.Item(2, 10).Value = 0.4216474758156494
$ws.Cells.Item(2, 11).Value = -0.2871398212255388

$ws.Cells.Item(3, 2).Value = -9.334904851338155
$ws.Cells.Item(3, 3).Value = -1.178654927766538
$ws.Cells.Item(3, 4).Value = 0.4008618608037394
$ws.Cells.Item(3, 5).Value = -1.854911900631668
$ws.Cells.Item(3, 6).Value = -1.131791633425802
$ws.Cells.Item(3, 7).Value = -0.333149228892605
$ws.Cells.Item(3, 8).Value = -0.4111884309045327
$ws.Cells.Item(3, 9).Value = 0.3169997218817192
$ws.Cells.Item(3, 10).Value = -0.3883310997898122
$ws.Cells.Item(3, 11).Value = 0.1055803449480684

$ws.Cells.Item(4, 2).Value = -8.343884810353487
$ws.Cells.Item(4, 3).Value = -5.373792936743199
$ws.Cells.Item(4, 4).Value = -6.496969861563686
$ws.Cells.Item(4, 5).Value = -4.882319462057946
$ws.Cells.Item(4, 6).Value = -3.389234397706988
$ws.Cells.Item(4, 7).Value = -2.930893202069679
$ws.Cells.Item(4, 8).Value = -1.79047838906512
$ws.Cells.Item(4, 9).Value = -2.180036841707216
$ws.Cells.Item(4, 10).Value = -1.444748162954129
$ws.Cells.Item(4, 11).Value = -0.7897867854766354

$ws.Cells.Item(5, 2).Value = -3.693611790173975
$ws.Cells.Item(5, 3).Value = 0.9190749009450181
$ws.Cells.Item(5, 4).Value = -2.634263145832853
$ws.Cells.Item(5, 5).Value = 0.5921266174821862
$ws.Cells.Item(5, 6).Value = -0.9527113001597033
$ws.Cells.Item(5, 7).Value = 0.6215278570052007
$ws.Cells.Item(5, 8).Value = -0.5900264632062522
$ws.Cells.Item(5, 9).Value = 0.1978243833556703
$ws.Cells.Item(5, 10).Value = 0.4935358542306867
$ws.Cells.Item(5, 11).Value = 0.9088263634637752

$ws.Cells.Item(6, 2).Value = -2.808131625772977
$ws.Cells.Item(6, 3).Value = -0.77139249106923
$ws.Cells.Item(6, 4).Value = -0.1687823066887937
$ws.Cells.Item(6, 5).Value = -0.4024803045374513
$ws.Cells.Item(6, 6).Value = 0.3742746160819611
$ws.Cells.Item(6, 7).Value = -0.3294423225852309
$ws.Cells.Item(6, 8).Value = 0.1540433631259821
$ws.Cells.Item(6, 9).Value = 0.6252233516998879
$ws.Cells.Item(6, 10).Value = 0.9388562853074672
$ws.Cells.Item(6, 11).Value = 0.0600786799400978

$ws.Cells.Item(7, 2).Value = -0.3631600374474599
$ws.Cells.Item(7, 3).Value = 0.08284073719760038
$ws.Cells.Item(7, 4).Value = -0.4843863583605866
$ws.Cells.Item(7, 5).Value = 0.4493666654632467
$ws.Cells.Item(7, 6).Value = -0.2235351620715071
$ws.Cells.Item(7, 7).Value = 0.2124600092846375
$ws.Cells.Item(7, 8).Value = 0.6916620355299794
$ws.Cells.Item(7, 9).Value = 1.013538243918304
$ws.Cells.Item(7, 10).Value = 0.1303919718410766
$ws.Cells.Item(7, 11).Value = -0.5999993966742521

$ws.Cells.Item(8, 2).Value = -0.02265847465121479
$ws.Cells.Item(8, 3).Value = -0.6104916423259444
$ws.Cells.Item(8, 4).Value = 0.6176147340964784
$ws.Cells.Item(8, 5).Value = -0.1628710694676248
$ws.Cells.Item(8, 6).Value = 0.2235928486592704
$ws.Cells.Item(8, 7).Value = 0.75452478685691
$ws.Cells.Item(8, 8).Value = 1.070889771712047
$ws.Cells.Item(8, 9).Value = 0.1748971488454556
$ws.Cells.Item(8, 10).Value = -0.5489255256478335
$ws.Cells.Item(8, 11).Value = 1.115476138671014

$ws.Cells.Item(9, 2).Value = -1.214932262173683
$ws.Cells.Item(9, 3).Value = 0.553411356741785
$ws.Cells.Item(9, 4).Value = 0.128676238589592
$ws.Cells.Item(9, 5).Value = 0.1858859725733368
$ws.Cells.Item(9, 6).Value = 0.7459407667834489
$ws.Cells.Item(9, 7).Value = 1.154630849309415
$ws.Cells.Item(9, 8).Value = 0.2097242937735006
$ws.Cells.Item(9, 9).Value = -0.5223611425955365
$ws.Cells.Item(9, 10).Value = 1.160974979417613
$ws.Cells.Item(9, 11).Value = 0.5338127914751168

$ws.Cells.Item(10, 2).Value = 0.4657705279028737
$ws.Cells.Item(10, 3).Value = 0.07625703655992344
$ws.Cells.Item(10, 4).Value = 0.2445323695583424
$ws.Cells.Item(10, 5).Value = 0.7508699121205502
$ws.Cells.Item(10, 6).Value = 1.14353756626052
$ws.Cells.Item(10, 7).Value = 0.2224909929985871
$ws.Cells.Item(10, 8).Value = -0.5136315835481774
$ws.Cells.Item(10, 9).Value = 1.163476923566464
$ws.Cells.Item(10, 10).Value = 0.540088728946798
$ws.Cells.Item(10, 11).Value = 0.7747069251866952

$ws.Cells.Item(11, 2).Value = 0.09868501619293202
$ws.Cells.Item(11, 3).Value = 0.2471359091207911
$ws.Cells.Item(11, 4).Value = 0.728802226466535
$ws.Cells.Item(11, 5).Value = 1.137555585956165
$ws.Cells.Item(11, 6).Value = 0.217147742816176
$ws.Cells.Item(11, 7).Value = -0.5240001169185436
$ws.Cells.Item(11, 8).Value = 1.154999601442841
$ws.Cells.Item(11, 9).Value = 0.5323164206399872
$ws.Cells.Item(11, 10).Value = 0.7660829516537266
$ws.Cells.Item(11, 11).Value = 0.6652412546050546

$ws.Cells.Item(12, 2).Value = 0.2685604805780112
$ws.Cells.Item(12, 3).Value = 0.8489800752657474
$ws.Cells.Item(12, 4).Value = 1.048361241542382
$ws.Cells.Item(12, 5).Value = 0.1874207389675982
$ws.Cells.Item(12, 6).Value = -0.5153739968886335
$ws.Cells.Item(12, 7).Value = 1.129059428017406
$ws.Cells.Item(12, 8).Value = 0.5095378576117597
$ws.Cells.Item(12, 9).Value = 0.7524649165051087
$ws.Cells.Item(12, 10).Value = 0.6467204574198788
$ws.Cells.Item(12, 11).Value = -0.1454929044188731

$ws.Cells.Item(13, 2).Value = 0.806659442945358
$ws.Cells.Item(13, 3).Value = 1.019156767686649
$ws.Cells.Item(13, 4).Value = 0.184400139162677
$ws.Cells.Item(13, 5).Value = -0.534235280500114
$ws.Cells.Item(13, 6).Value = 1.108728296952201
$ws.Cells.Item(13, 7).Value = 0.4950077021775119
$ws.Cells.Item(13, 8).Value = 0.7358238710128409
$ws.Cells.Item(13, 9).Value = 0.6291366119911286
$ws.Cells.Item(13, 10).Value = -0.1619667879489159
$ws.Cells.Item(13, 11).Value = 0.449372724506711

$ws.Cells.Item(14, 2).Value = 1.359354508304559
$ws.Cells.Item(14, 3).Value = 0.2597712009466141
$ws.Cells.Item(14, 4).Value = -0.7253285668131131
$ws.Cells.Item(14, 5).Value = 1.133422209538934
$ws.Cells.Item(14, 6).Value = 0.5043871057338079
$ws.Cells.Item(14, 7).Value = 0.6810130724014498
$ws.Cells.Item(14, 8).Value = 0.6102425211024588
$ws.Cells.Item(14, 9).Value = -0.1768307887639616
$ws.Cells.Item(14, 10).Value = 0.4206610138773402
$ws.Cells.Item(14, 11).Value = 0.2638965897873631

$ws.Cells.Item(15, 2).Value = 0.7104660729368646
$ws.Cells.Item(15, 3).Value = -0.6782178646007546
$ws.Cells.Item(15, 4).Value = 0.8931389787481262
$ws.Cells.Item(15, 5).Value = 0.5403321867996707
$ws.Cells.Item(15, 6).Value = 0.6789623149618569
$ws.Cells.Item(15, 7).Value = 0.5336742816664286
$ws.Cells.Item(15, 8).Value = -0.204230183139245
$ws.Cells.Item(15, 9).Value = 0.3943719636796149
$ws.Cells.Item(15, 10).Value = 0.2204992990740305
$ws.Cells.Item(15, 11).ClearContents()

$ws.Cells.Item(16, 2).Value = -0.3658922776772162
$ws.Cells.Item(16, 3).Value = 1.026427718483651
$ws.Cells.Item(16, 4).Value = 0.3612537649834791
$ws.Cells.Item(16, 5).Value = 0.707156553732424
$ws.Cells.Item(16, 6).Value = 0.5696995928588221
$ws.Cells.Item(16, 7).Value = -0.2417717533434586
$ws.Cells.Item(16, 8).Value = 0.3888296245922537
$ws.Cells.Item(16, 9).Value = 0.2246746280127792
$ws.Cells.Item(16, 10).ClearContents()
$ws.Cells.Item(16, 11).ClearContents()

$ws.Cells.Item(17, 2).Value = 1.262018209591492
$ws.Cells.Item(17, 3).Value = 0.4459325058577887
$ws.Cells.Item(17, 4).Value = 0.5615561870987069
$ws.Cells.Item(17, 5).Value = 0.5796533357180647
$ws.Cells.Item(17, 6).Value = -0.2264065836439137
$ws.Cells.Item(17, 7).Value = 0.3499812896348306
$ws.Cells.Item(17, 8).Value = 0.2088288189855932
$ws.Cells.Item(17, 9).ClearContents()
$ws.Cells.Item(17, 10).ClearContents()
$ws.Cells.Item(17, 11).ClearContents()

$ws.Cells.Item(18, 2).Value = 0.7569566923391715
$ws.Cells.Item(18, 3).Value = 0.6786449615099022
$ws.Cells.Item(18, 4).Value = 0.4171891942684979
$ws.Cells.Item(18, 5).Value = -0.1974476331787121
$ws.Cells.Item(18, 6).Value = 0.3857241620897341
$ws.Cells.Item(18, 7).Value = 0.1775011726019661
$ws.Cells.Item(18, 8).ClearContents()
$ws.Cells.Item(18, 9).ClearContents()
$ws.Cells.Item(18, 10).ClearContents()
$ws.Cells.Item(18, 11).ClearContents()

$ws.Cells.Item(19, 2).Value = 0.9254701389140165
$ws.Cells.Item(19, 3).Value = 0.4348450618063874
$ws.Cells.Item(19, 4).Value = -0.2920698722897066
$ws.Cells.Item(19, 5).Value = 0.4179415503382142
$ws.Cells.Item(19, 6).Value = 0.1898892984296834
$ws.Cells.Item(19, 7).ClearContents()
$ws.Cells.Item(19, 8).ClearContents()
$ws.Cells.Item(19, 9).ClearContents()
$ws.Cells.Item(19, 10).ClearContents()
$ws.Cells.Item(19, 11).ClearContents()

$ws.Cells.Item(20, 2).Value = 0.6745214212225993
$ws.Cells.Item(20, 3).Value = -0.2070456288204931
$ws.Cells.Item(20, 4).Value = 0.3014569719802002
$ws.Cells.Item(20, 5).Value = 0.2049945700815359
$ws.Cells.Item(20, 6).ClearContents()
$ws.Cells.Item(20, 7).ClearContents()
$ws.Cells.Item(20, 8).ClearContents()
$ws.Cells.Item(20, 9).ClearContents()
$ws.Cells.Item(20, 10).ClearContents()
$ws.Cells.Item(20, 11).ClearContents()

$ws.Cells.Item(21, 2).Value = -0.04218555178640582
$ws.Cells.Item(21, 3).Value = 0.3149942442281164
$ws.Cells.Item(21, 4).Value = 0.1420216510915729
$ws.Cells.Item(21, 5).ClearContents()
$ws.Cells.Item(21, 6).ClearContents()
$ws.Cells.Item(21, 7).ClearContents()
$ws.Cells.Item(21, 8).ClearContents()
$ws.Cells.Item(21, 9).ClearContents()
$ws.Cells.Item(21, 10).ClearContents()
$ws.Cells.Item(21, 11).ClearContents()

$ws.Cells.Item(22, 2).Value = 0.5688432860935244
$ws.Cells.Item(22, 3).Value = 0.2413397012736094
$ws.Cells.Item(22, 4).ClearContents()
$ws.Cells.Item(22, 5).ClearContents()
$ws.Cells.Item(22, 6).ClearContents()
$ws.Cells.Item(22, 7).ClearContents()
$ws.Cells.Item(22, 8).ClearContents()
$ws.Cells.Item(22, 9).ClearContents()
$ws.Cells.Item(22, 10).ClearContents()
$ws.Cells.Item(22, 11).ClearContents()

$ws.Cells.Item(23, 2).Value = 0.2853993925130583
$ws.Cells.Item(23, 3).ClearContents()
$ws.Cells.Item(23, 4).ClearContents()
$ws.Cells.Item(23, 5).ClearContents()
$ws.Cells.Item(23, 6).ClearContents()
$ws.Cells.Item(23, 7).ClearContents()
$ws.Cells.Item(23, 8).ClearContents()
$ws.Cells.Item(23, 9).ClearContents()
$ws.Cells.Item(23, 10).ClearContents()
$ws.Cells.Item(23, 11).ClearContents()

$ws.Cells.Item(24, 2).ClearContents()
$ws.Cells.Item(24, 3).ClearContents()
$ws.Cells.Item(24, 4).ClearContents()
$ws.Cells.Item(24, 5).ClearContents()
$ws.Cells.Item(24, 6).ClearContents()
$ws.Cells.Item(24, 7).ClearContents()
$ws.Cells.Item(24, 8).ClearContents()
$ws.Cells.Item(24, 9).ClearContents()
$ws.Cells.Item(24, 10).ClearContents()
$ws.Cells.Item(24, 11).ClearContents()
